$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shared-string (rich text) updates
# ---------------------------------------------------------------------------

# A8: "Volume 31   Number  31" -> "Volume 31   Number  32"
$cellA8 = $ws.Range("A8")
$cellA8.Value = "Volume 31   Number  32"
$cellA8.Characters(1, 7).Font.Size = 10
$cellA8.Characters(1, 7).Font.Name = "Andale WT"
$cellA8.Characters(8, 2).Font.Size = 10
$cellA8.Characters(8, 2).Font.Name = "Andale WT"
$cellA8.Characters(10, 11).Font.Size = 10
$cellA8.Characters(10, 11).Font.Name = "Andale WT"
$cellA8.Characters(21, 2).Font.Size = 10
$cellA8.Characters(21, 2).Font.Name = "Andale WT"

# C9: "Report Covering the Week  7/29/2024  Through  8/4/2024"
#  -> "Report Covering the Week  8/5/2024  Through  8/11/2024"
$cellC9 = $ws.Range("C9")
$cellC9.Value = "Report Covering the Week  8/5/2024  Through  8/11/2024"
$cellC9.Characters(1, 26).Font.Size = 10
$cellC9.Characters(1, 26).Font.Name = "Andale WT"
$cellC9.Characters(27, 8).Font.Size = 10
$cellC9.Characters(27, 8).Font.Name = "Andale WT"
$cellC9.Characters(35, 11).Font.Size = 10
$cellC9.Characters(35, 11).Font.Name = "Andale WT"
$cellC9.Characters(46, 9).Font.Size = 10
$cellC9.Characters(46, 9).Font.Name = "Andale WT"

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Range("L15").Value = -22.222222222222

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 63
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = -10
$ws.Range("L16").Value = -25.882352941176
$ws.Range("M16").Value = -22.222222222222
$ws.Range("N16").Value = -82.83378746594

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = -31.25
$ws.Range("I17").Value = 79
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = 12.857142857142
$ws.Range("L17").Value = -17.708333333333
$ws.Range("M17").Value = 61.224489795918
$ws.Range("N17").Value = -50.314465408805

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -46.666666666666
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = -25.423728813559
$ws.Range("L18").Value = -38.888888888888
$ws.Range("M18").Value = -42.857142857142
$ws.Range("N18").Value = -87.374461979913

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 18.367346938775
$ws.Range("I19").Value = 394
$ws.Range("J19").Value = 411
$ws.Range("K19").Value = -4.136253041362
$ws.Range("L19").Value = 10.674157303370
$ws.Range("M19").Value = 137.349397590361
$ws.Range("N19").Value = 81.566820276497

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -66.666666666666
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -61.111111111111
$ws.Range("I20").Value = 58
$ws.Range("J20").Value = 104
$ws.Range("K20").Value = -44.230769230769
$ws.Range("L20").Value = -45.794392523364
$ws.Range("M20").Value = -35.555555555555
$ws.Range("N20").Value = -89.679715302491

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -25.925925925925
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = -11.926605504587
$ws.Range("I21").Value = 689
$ws.Range("J21").Value = 776
$ws.Range("K21").Value = -11.211340206185
$ws.Range("L21").Value = -13.659147869674
$ws.Range("M21").Value = 27.356746765249
$ws.Range("N21").Value = -65.789473684210

# ---------------------------------------------------------------------------
# Row 23 (type changes on C23, D23, E23)
# ---------------------------------------------------------------------------
# C23: numeric -> text "0" (style 14, shared string "0")
$ws.Range("D14").Copy($ws.Range("C23"))
# D23: text "0" -> numeric 2 (style 16)
$ws.Range("C16").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 2
# E23: text "***.*" -> numeric -100 (style 15)
$ws.Range("E16").Copy($ws.Range("E23"))
$ws.Range("E23").Value = -100

$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = 5.882352941176
$ws.Range("L23").Value = 20
$ws.Range("M23").Value = -5.263157894736

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 23.529411764705
$ws.Range("I24").Value = 606
$ws.Range("J24").Value = 564
$ws.Range("K24").Value = 7.446808510638
$ws.Range("L24").Value = -3.04
$ws.Range("M24").Value = 73.142857142857

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 18
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 107.142857142857
$ws.Range("I25").Value = 358
$ws.Range("J25").Value = 289
$ws.Range("K25").Value = 23.875432525951
$ws.Range("L25").Value = -0.555555555555

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -77.777777777777
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = -35
$ws.Range("I26").Value = 165
$ws.Range("J26").Value = 148
$ws.Range("K26").Value = 11.486486486486
$ws.Range("L26").Value = -4.624277456647
$ws.Range("M26").Value = 18.705035971223

# ---------------------------------------------------------------------------
# Row 27 (type changes on D27, E27)
# ---------------------------------------------------------------------------
# D27: numeric -> text "0" (style 14, shared string "0")
$ws.Range("D14").Copy($ws.Range("D27"))
# E27: numeric -> text "***.*" (style 14, shared string "***.*")
$ws.Range("E14").Copy($ws.Range("E27"))

$ws.Range("L27").Value = -27.272727272727

# ---------------------------------------------------------------------------
# Row 28 (type changes on C28, D28, E28)
# ---------------------------------------------------------------------------
# C28: numeric -> text "0" (style 14, shared string "0")
$ws.Range("D14").Copy($ws.Range("C28"))
# D28: text "0" -> numeric 1 (style 16)
$ws.Range("C16").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
# E28: text "***.*" -> numeric -100 (style 15)
$ws.Range("E16").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100

$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = -25

# ---------------------------------------------------------------------------
# Row 29
# ---------------------------------------------------------------------------
$ws.Range("L29").Value = -80
$ws.Range("N29").Value = -90

# ---------------------------------------------------------------------------
# Row 30
# ---------------------------------------------------------------------------
$ws.Range("L30").Value = -80
$ws.Range("N30").Value = -87.5
